# Generate Report for Handoff
# Marks the 6869bf9d...md file as "Ready for handoff" (it is now up to date
# for hand-off instead of still showing the prior "Handed back" status), and
# records the new handoff timestamps + an out-of-date handback warning on
# the per-locale detail sheets.

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5a069f12556b750676ac0f711c289b3dd1eb4330/e2e/6869bf9d-595d-4cb8-8a59-5bf12cf690da.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a203e470079b6290647fb368c6293a1c35765cc2/e2e/6869bf9d-595d-4cb8-8a59-5bf12cf690da.md."

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $readyForHandoff
$overview.Range("F3").Value = $readyForHandoff
$overview.Range("G3").Value = "2016-09-02 20:56:20"

# --- zh-cn detail sheet ----------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $readyForHandoff
$zhcn.Range("H2").Value = "2016-09-02 20:56:15"
$zhcn.Range("C3").Value = $readyForHandoff
$zhcn.Range("H3").Value = "2016-09-02 20:56:15"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns(16).ColumnWidth = 39.16666666666667

# --- de-de detail sheet ----------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $readyForHandoff
$dede.Range("H3").Value = "2016-09-02 20:56:20"
$dede.Range("P3").Value = $errorDetail
$dede.Columns(16).ColumnWidth = 39.16666666666667
